$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column BO (nomor_s) values from 2010 to 2338 for data rows 6 through 196
$ws.Range("BO6:BO196").Value = 2338
